{"js": "// Apply the four text replacements described by the diff using the\n// Word JavaScript API (Office.js). `context` is the Word.RequestContext.\n\nconst replacements = [\n  {\n    old: \"Ativa\u00e7\u00e3o: 01/01/2018\",\n    new: \"Ativa\u00e7\u00e3o: 01/01/2021\",\n  },\n  {\n    old:\n      \"1) Apresentar as 1\u00aa e 2\u00aa Leis da Termodin\u00e2mica e aplic\u00e1-las a problemas reais de engenharia;\" +\n      \"2) Calcular ciclos t\u00e9rmicos, ciclos de refrigera\u00e7\u00e3o e combust\u00e3o, para que o Engenheiro de \" +\n      \"Materiais possa otimizar a efici\u00eancia de ciclos t\u00e9rmicos usando materiais que se adequem \u00e1s \" +\n      \"condi\u00e7\u00f5es de projeto dos ciclos.\",\n    new:\n      \"Abordar os princ\u00edpios b\u00e1sicos da termodin\u00e2mica de forma que os estudantes e futuros \" +\n      \"engenheiros tenham um entendimento claro e s\u00f3lido sobre estes princ\u00edpios. Apresentar \" +\n      \"diversos exemplos de engenharia do mundo real e de como a termodin\u00e2mica \u00e9 aplicada na \" +\n      \"pr\u00e1tica de engenharia. Enfatizar a compreens\u00e3o da termodin\u00e2mica baseada na F\u00edsica e em \" +\n      \"argumentos f\u00edsicos, buscando incentivar o entendimento mais profundo da termodin\u00e2mica.\",\n  },\n  {\n    old:\n      \"1. Conceitos, Defini\u00e7\u00f5es e Propriedades de uma subst\u00e2ncia pura2. Trabalho e Calor3. 1\u00aa Lei \" +\n      \"de Termodin\u00e2mica4. 2\u00aa Lei da Termodin\u00e2mica5. Entropia6. Ciclo Motores e de Refrigera\u00e7\u00e3o7. \" +\n      \"Projeto sobre Gera\u00e7\u00e3o de Energia Termoel\u00e9trica: Ciclo Simples (vapor), Ciclo Combinado \" +\n      \"(turbina a g\u00e1s/caldeira - turbina a vapor), Ciclos de refrigera\u00e7\u00e3o e de gera\u00e7\u00e3o de pot\u00eancia \" +\n      \"combinados.\",\n    new:\n      \"1. Termodin\u00e2mica e Energia. 2. Import\u00e2ncia das unidades e an\u00e1lise dimensional.3. Sistemas e \" +\n      \"volumes de controle. 4. Equipamentos dom\u00e9sticos e a Termodin\u00e2mica. 5. Propriedades de um \" +\n      \"sistema: estados termodin\u00e2micos e equil\u00edbrio. 6. Efici\u00eancia na convers\u00e3o de energia. 7. \" +\n      \"Processos e ciclos t\u00e9rmicos. 8. Termodin\u00e2mica e o meio ambiente.\",\n  },\n  {\n    old:\n      \"1.Conceitos, defini\u00e7\u00f5es e propriedades de uma subst\u00e2ncia pura; 2.Trabalho e calor;  3.1\u00aa Lei \" +\n      \"da termodin\u00e2mica: Teoria e aplica\u00e7\u00e3o a volumes de controle; 4.2\u00aa Lei da termodin\u00e2mica: \" +\n      \"Entropia5.2\u00aa Lei da termodin\u00e2mica: Aplica\u00e7\u00e3o a volumes de controle;6.Ciclos motores Ciclos \" +\n      \"de refrigera\u00e7\u00e3o;7.Projeto sobre ciclo simples: Vapor; Projeto sobre ciclos combinados: \" +\n      \"Turbina a g\u00e1s, turbina a vapor, Ciclos de refrigera\u00e7\u00e3o e de gera\u00e7\u00e3o de pot\u00eancia combinados\",\n    new:\n      \"1. Termodin\u00e2mica e Energia: formas de energia e transfer\u00eancia de energia por calor e \" +\n      \"trabalho; formas mec\u00e2nicas de trabalho. 2. Sistema de Unidades e An\u00e1lise Dimensional: \" +\n      \"import\u00e2ncia na engenharia de m\u00e1quinas. 3. Sistemas e volumes de controle: dispositivos \" +\n      \"ativos e passivos. 4. Propriedades de um sistema. Estados e equil\u00edbrio: diagramas de \" +\n      \"propriedades para processos com mudan\u00e7a de fase; equil\u00edbrio de estado do g\u00e1s ideal; fator \" +\n      \"de compressibilidade; press\u00e3o de vapor e press\u00e3o de equil\u00edbrio; calores espec\u00edficos. 5. \" +\n      \"Balan\u00e7o de energia em sistemas fechados e em volumes de controle: trabalho de fluxo e \" +\n      \"energia de escoamento de um fluido; regime permanente e transiente. 6. M\u00e1quinas t\u00e9rmicas e \" +\n      \"refrigeradores e a 2\u00aa. Lei da Termodin\u00e2mica: princ\u00edpios e ciclos de Carnot; entropia e \" +\n      \"varia\u00e7\u00e3o de entropia em s\u00f3lidos, l\u00edquidos e gases. 7. Efici\u00eancia na convers\u00e3o de energia. \" +\n      \"Efici\u00eancia t\u00e9rmica. Efici\u00eancia de m\u00e1quinas. Efici\u00eancia isoentr\u00f3pica em dispositivos com \" +\n      \"escoamento em regime permanente. Balan\u00e7o de entropia. 8. Processo e ciclos: Ciclos de \" +\n      \"pot\u00eancia a g\u00e1s: Otto, Diesel, Stirling, Ericsson, Brayton e suas varia\u00e7\u00f5es. Ciclos de \" +\n      \"pot\u00eancia a vapor e ciclos combinados g\u00e1s-vapor: Rankine ideal; afastamento da condi\u00e7\u00e3o \" +\n      \"ideal; efici\u00eancia do ciclo Rankine com e sem modifica\u00e7\u00f5es; cogera\u00e7\u00e3o. Ciclos de \" +\n      \"refrigera\u00e7\u00e3o e sistemas de bombas de calor: sistemas a g\u00e1s e por absor\u00e7\u00e3o. 9. Economia de \" +\n      \"energia: benef\u00edcios ao meio ambiente.\",\n  },\n];\n\nfor (const { old, new: replacement } of replacements) {\n  const results = context.document.body.search(old, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${old.slice(0, 40)}...`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replacement, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the four text replacements described by the diff using the Word\n# COM object model. `$word` is the Word.Application and `$word.ActiveDocument`\n# (also available as `$d` / `$doc`) is the open document.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        Old = \"Ativa\u00e7\u00e3o: 01/01/2018\"\n        New = \"Ativa\u00e7\u00e3o: 01/01/2021\"\n    },\n    @{\n        Old = \"1) Apresentar as 1\u00aa e 2\u00aa Leis da Termodin\u00e2mica e aplic\u00e1-las a problemas reais de engenharia;2) Calcular ciclos t\u00e9rmicos, ciclos de refrigera\u00e7\u00e3o e combust\u00e3o, para que o Engenheiro de Materiais possa otimizar a efici\u00eancia de ciclos t\u00e9rmicos usando materiais que se adequem \u00e1s condi\u00e7\u00f5es de projeto dos ciclos.\"\n        New = \"Abordar os princ\u00edpios b\u00e1sicos da termodin\u00e2mica de forma que os estudantes e futuros engenheiros tenham um entendimento claro e s\u00f3lido sobre estes princ\u00edpios. Apresentar diversos exemplos de engenharia do mundo real e de como a termodin\u00e2mica \u00e9 aplicada na pr\u00e1tica de engenharia. Enfatizar a compreens\u00e3o da termodin\u00e2mica baseada na F\u00edsica e em argumentos f\u00edsicos, buscando incentivar o entendimento mais profundo da termodin\u00e2mica.\"\n    },\n    @{\n        Old = \"1. Conceitos, Defini\u00e7\u00f5es e Propriedades de uma subst\u00e2ncia pura2. Trabalho e Calor3. 1\u00aa Lei de Termodin\u00e2mica4. 2\u00aa Lei da Termodin\u00e2mica5. Entropia6. Ciclo Motores e de Refrigera\u00e7\u00e3o7. Projeto sobre Gera\u00e7\u00e3o de Energia Termoel\u00e9trica: Ciclo Simples (vapor), Ciclo Combinado (turbina a g\u00e1s/caldeira - turbina a vapor), Ciclos de refrigera\u00e7\u00e3o e de gera\u00e7\u00e3o de pot\u00eancia combinados.\"\n        New = \"1. Termodin\u00e2mica e Energia. 2. Import\u00e2ncia das unidades e an\u00e1lise dimensional.3. Sistemas e volumes de controle. 4. Equipamentos dom\u00e9sticos e a Termodin\u00e2mica. 5. Propriedades de um sistema: estados termodin\u00e2micos e equil\u00edbrio. 6. Efici\u00eancia na convers\u00e3o de energia. 7. Processos e ciclos t\u00e9rmicos. 8. Termodin\u00e2mica e o meio ambiente.\"\n    },\n    @{\n        Old = \"1.Conceitos, defini\u00e7\u00f5es e propriedades de uma subst\u00e2ncia pura; 2.Trabalho e calor;  3.1\u00aa Lei da termodin\u00e2mica: Teoria e aplica\u00e7\u00e3o a volumes de controle; 4.2\u00aa Lei da termodin\u00e2mica: Entropia5.2\u00aa Lei da termodin\u00e2mica: Aplica\u00e7\u00e3o a volumes de controle;6.Ciclos motores Ciclos de refrigera\u00e7\u00e3o;7.Projeto sobre ciclo simples: Vapor; Projeto sobre ciclos combinados: Turbina a g\u00e1s, turbina a vapor, Ciclos de refrigera\u00e7\u00e3o e de gera\u00e7\u00e3o de pot\u00eancia combinados\"\n        New = \"1. Termodin\u00e2mica e Energia: formas de energia e transfer\u00eancia de energia por calor e trabalho; formas mec\u00e2nicas de trabalho. 2. Sistema de Unidades e An\u00e1lise Dimensional: import\u00e2ncia na engenharia de m\u00e1quinas. 3. Sistemas e volumes de controle: dispositivos ativos e passivos. 4. Propriedades de um sistema. Estados e equil\u00edbrio: diagramas de propriedades para processos com mudan\u00e7a de fase; equil\u00edbrio de estado do g\u00e1s ideal; fator de compressibilidade; press\u00e3o de vapor e press\u00e3o de equil\u00edbrio; calores espec\u00edficos. 5. Balan\u00e7o de energia em sistemas fechados e em volumes de controle: trabalho de fluxo e energia de escoamento de um fluido; regime permanente e transiente. 6. M\u00e1quinas t\u00e9rmicas e refrigeradores e a 2\u00aa. Lei da Termodin\u00e2mica: princ\u00edpios e ciclos de Carnot; entropia e varia\u00e7\u00e3o de entropia em s\u00f3lidos, l\u00edquidos e gases. 7. Efici\u00eancia na convers\u00e3o de energia. Efici\u00eancia t\u00e9rmica. Efici\u00eancia de m\u00e1quinas. Efici\u00eancia isoentr\u00f3pica em dispositivos com escoamento em regime permanente. Balan\u00e7o de entropia. 8. Processo e ciclos: Ciclos de pot\u00eancia a g\u00e1s: Otto, Diesel, Stirling, Ericsson, Brayton e suas varia\u00e7\u00f5es. Ciclos de pot\u00eancia a vapor e ciclos combinados g\u00e1s-vapor: Rankine ideal; afastamento da condi\u00e7\u00e3o ideal; efici\u00eancia do ciclo Rankine com e sem modifica\u00e7\u00f5es; cogera\u00e7\u00e3o. Ciclos de refrigera\u00e7\u00e3o e sistemas de bombas de calor: sistemas a g\u00e1s e por absor\u00e7\u00e3o. 9. Economia de energia: benef\u00edcios ao meio ambiente.\"\n    }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $($r.Old.Substring(0, 40))\"\n    }\n}\n"}
